# [Fonds de solidarite] Add 2022-06-09 data
# Updates nombre_aides (column C) and montant_total (column E) values
# for the rows whose underlying data changed with the new extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 53;  C = 141681;  E = 590068269 },
    @{ Row = 83;  C = 3414;    E = 115758395 },
    @{ Row = 91;  C = 151142;  E = 482442285 },
    @{ Row = 92;  C = 409155;  E = 1595512102 },
    @{ Row = 95;  C = 50776;   E = 932943241 },
    @{ Row = 96;  C = 17297;   E = 794564257 },
    @{ Row = 104; C = 135245;  E = 272235386 },
    @{ Row = 116; C = 4563;    E = 20627806 },
    @{ Row = 174; C = 226098;  E = 900654909 },
    @{ Row = 177; C = 14719;  E = 251626957 }
)

foreach ($u in $updates) {
    $ws.Range("C" + $u.Row).Value = $u.C
    $ws.Range("E" + $u.Row).Value = $u.E
}
